# Efetivando as SARs aprovadas no EVE-EOR.
# Rename the status text "Efetivada e Encerrada" -> "Efetivada no EOR" and
# apply this new status to the three SAR rows on "Histórico de SARs" that
# previously read "Aprovada para Resolução". Also update the last
# selected cell on both sheets, matching the target workbook.

$wb = $excel.ActiveWorkbook

$wsSARs = $wb.Worksheets.Item("Histórico de SARs")
$wsLegenda = $wb.Worksheets.Item("Legenda do Documento")

# Update the legend entry text itself (shared string reused by C5 on the
# legend sheet as well as the status column on the SARs sheet).
$wsLegenda.Range("C5").Value = "Efetivada no EOR"

# Apply the renamed status to the approved SARs.
$wsSARs.Range("C2").Value = "Efetivada no EOR"
$wsSARs.Range("C3").Value = "Efetivada no EOR"
$wsSARs.Range("C4").Value = "Efetivada no EOR"

# Update the last active selection on each sheet to match the new state.
$wsSARs.Range("C10").Select()
$wsLegenda.Range("C5").Select()

$wsSARs.Activate()
